$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.224.93"
$ws.Cells.Item(2, 5).Value = "  +0.06%  "

$ws.Cells.Item(3, 4).Value = "1.902.45"
$ws.Cells.Item(3, 5).Value = "  +0.48%  "

$ws.Cells.Item(4, 5).Value = "  -0.22%  "

$ws.Cells.Item(5, 4).Value = "'306.33"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.37%  "

$ws.Cells.Item(6, 5).Value = "  -0.05%  "

$ws.Cells.Item(7, 4).Value = "'0.5369"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +3.38%  "

$ws.Cells.Item(8, 4).Value = "'0.3811"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +1.48%  "

$ws.Cells.Item(9, 4).Value = "'0.07275"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.03%  "

$ws.Cells.Item(10, 4).Value = "'22.19"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +4.83%  "

$ws.Cells.Item(11, 4).Value = "'0.9031"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.41%  "

$ws.Cells.Item(12, 4).Value = "'0.08187"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.72%  "

$ws.Cells.Item(13, 4).Value = "'95.96"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.36%  "

$ws.Cells.Item(14, 4).Value = "'5.339"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +1.05%  "

$ws.Cells.Item(15, 4).Value = "'0.9997"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.40%  "

$ws.Cells.Item(16, 4).Value = "'14.84"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.05%  "

$ws.Cells.Item(17, 4).Value = "'0.000008643"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.63%  "

$ws.Cells.Item(18, 5).Value = "  -0.04%  "

$ws.Cells.Item(19, 4).Value = "27.256.48"
$ws.Cells.Item(19, 5).Value = "  +0.07%  "

$ws.Cells.Item(20, 5).Value = "  -0.92%  "

$ws.Cells.Item(21, 4).Value = "1.098.23"
$ws.Cells.Item(21, 5).Value = "  -42.03%  "

$ws.Cells.Item(22, 5).Value = "  +0.87%  "

$ws.Cells.Item(23, 4).Value = "'6.494"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.60%  "

$ws.Cells.Item(24, 4).Value = "'149.38"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +1.61%  "

$ws.Cells.Item(25, 4).Value = "'2.289"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.46%  "

$ws.Cells.Item(26, 5).Value = "  +0.62%  "

$ws.Cells.Item(27, 5).Value = "  -0.09%  "

$ws.Cells.Item(28, 4).Value = "'116.67"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.28%  "

$ws.Cells.Item(29, 4).Value = "'4.813"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.03%  "

$ws.Cells.Item(30, 5).Value = "  -4.02%  "

$ws.Cells.Item(31, 4).Value = "'0.09216"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.06%  "

$ws.Cells.Item(32, 4).Value = "'0.8324"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +4.81%  "

$ws.Cells.Item(33, 4).Value = "'0.05066"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.64%  "

$ws.Cells.Item(34, 5).Value = "  -0.24%  "

$ws.Cells.Item(35, 4).Value = "'3.003"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +1.46%  "

$ws.Cells.Item(36, 4).Value = "'3.336"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -3.31%  "

$ws.Cells.Item(37, 4).Value = "'2.687"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +3.46%  "

$ws.Cells.Item(38, 4).Value = "'0.5836"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +3.02%  "

$ws.Cells.Item(39, 5).Value = "  +0.91%  "

$ws.Cells.Item(40, 5).Value = "  +0.22%  "

$ws.Cells.Item(41, 4).Value = "'9.297"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +3.92%  "

$ws.Cells.Item(42, 4).Value = "'6.606"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.82%  "

$ws.Cells.Item(43, 4).Value = "'116.93"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +1.28%  "

$ws.Cells.Item(44, 4).Value = "'0.1522"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.63%  "

$ws.Cells.Item(45, 4).Value = "'0.5005"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +3.13%  "

$ws.Cells.Item(46, 4).Value = "'1.001"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.06%  "

$ws.Cells.Item(47, 4).Value = "'10.07"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.29%  "

$ws.Cells.Item(48, 4).Value = "'1.637"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.22%  "

$ws.Cells.Item(49, 4).Value = "'38.34"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.62%  "

$ws.Cells.Item(50, 4).Value = "'0.06165"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +3.64%  "

$ws.Cells.Item(51, 5).Value = "  -0.11%  "
